$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix X20: was stored as text "21", should become a real number 21 ---
$ws.Range("X20").Value = 21

# --- New row 21 ---
$ws.Range("R21").Value = "QUOTE-20251213091212"
$ws.Range("S21").Value = "nabeel"
$ws.Range("T21").Value = "Karachi Port"
$ws.Range("U21").Value = "Bandar Abbas"
$ws.Range("V21").Value = "Food Item"
$ws.Range("X21").Value = 12
$ws.Range("Z21").Value = "2025-12-13 09:12:12"
$ws.Range("AA21").Value = "Dry"

# --- New row 22 ---
$ws.Range("R22").Value = "QUOTE-20251213092116"
$ws.Range("S22").Value = "nabeel"
$ws.Range("T22").Value = "Karachi Port"
$ws.Range("U22").Value = "Bandar Abbas Port"
$ws.Range("V22").Value = "Food Item"
$ws.Range("X22").Value = 123
$ws.Range("Z22").Value = "2025-12-13 09:21:16"
$ws.Range("AA22").Value = "Dry"

# --- New row 23 ---
$ws.Range("R23").Value = "QUOTE-20251213100629"
$ws.Range("S23").Value = "nabeel"
$ws.Range("T23").Value = "Karachi Port"
$ws.Range("U23").Value = "Bandar Abbas Port"
$ws.Range("V23").Value = "Food Item"
$ws.Range("X23").Value = 123
$ws.Range("Z23").Value = "2025-12-13 10:06:29"
$ws.Range("AA23").Value = "Dry"

# --- New row 24 ---
$ws.Range("R24").Value = "QUOTE-20251213100836"
$ws.Range("S24").Value = "nabeel"
$ws.Range("T24").Value = "Karachi Port"
$ws.Range("U24").Value = "Bandar Abbas Port"
$ws.Range("V24").Value = "General Cargo"
$ws.Range("X24").Value = 23
$ws.Range("Z24").Value = "2025-12-13 10:08:36"
$ws.Range("AA24").Value = "20ft Dry"

# --- New row 25 ---
$ws.Range("R25").Value = "QUOTE-20251213101827"
$ws.Range("S25").Value = "nabeel"
$ws.Range("T25").Value = "Karachi Port"
$ws.Range("U25").Value = "Bandar Abbas Port"
$ws.Range("V25").Value = "Food Item"
$ws.Range("X25").Value = 1312
$ws.Range("Z25").Value = "2025-12-13 10:18:27"
$ws.Range("AA25").Value = "Dry"

# --- New row 26 ---
$ws.Range("R26").Value = "QUOTE-20251213101852"
$ws.Range("S26").Value = "nabeel"
$ws.Range("T26").Value = "Karachi Port"
$ws.Range("U26").Value = "Bandar Abbas Port"
$ws.Range("V26").Value = "General Cargo"
$ws.Range("X26").Value = 12
$ws.Range("Z26").Value = "2025-12-13 10:18:52"
$ws.Range("AA26").Value = "Dry"

# --- New row 27 ---
$ws.Range("R27").Value = "QUOTE-20251213103013"
$ws.Range("S27").Value = "nabeel"
$ws.Range("T27").Value = "Karachi Port"
$ws.Range("U27").Value = "Bandar Abbas Port"
$ws.Range("V27").Value = "General Cargo"
# X27 keeps "12" as text (unlike the other new rows, which use real numbers)
$ws.Range("X27").NumberFormat = "@"
$ws.Range("X27").Value = "12"
$ws.Range("X27").Style = "Normal"
$ws.Range("Z27").Value = "2025-12-13 10:30:13"
$ws.Range("AA27").Value = "Dry"
